$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 to I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-23
$data = @(
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(7, 8),
    @(5, 7),
    @(5, 6),
    @(7, 8),
    @(8, 8),
    @(6, 6),
    @(5, 6),
    @(7, 8),
    @(4, 5),
    @(6, 6),
    @(5, 5),
    @(5, 6),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}

$wb.Save()
